$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: extend the bottom border row into column L (reuses existing style) ---
$ws.Range("J3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# --- Row 4: new "2021" header cell, bold Times New Roman + medium-bottom border ---
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Borders.Item(8).LineStyle = -4142
$ws.Range("L4").Value = 2021

# --- Row 5: new data cell, regular Times New Roman, no border ---
$ws.Range("K4").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Font.Bold = $false
$ws.Range("L5").Borders.Item(8).LineStyle = -4142
$ws.Range("L5").Borders.Item(9).LineStyle = -4142
$ws.Range("L5").Value = 2.3

# --- Row 6: new data cell, regular Times New Roman, medium-bottom border ---
$ws.Range("K4").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Font.Bold = $false
$ws.Range("L6").Borders.Item(8).LineStyle = -4142
$ws.Range("L6").Value = 1.3

# --- Update selection to match the recorded UI state after the edit ---
$ws.Range("O5").Select() | Out-Null
